$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "in_service" flags (column E) for rows 11-15 from FALSE to TRUE
$ws.Range("E11:E15").Value = $true
